# "modified 2005 2006 2007" — trim the trailing 4 summary rows (Totals for
# Foreign/U.S. points + confidential/non-confidential totals) off the 2005
# sheet, pull the print area in to match, and leave the selection sitting
# on the (now-blank) rows that used to separate the table from the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 249-252 held the four "Total ..." summary lines (shared strings
# #237-#240) that the commit drops. Deleting the whole rows shifts
# everything below (the blank style-only rows 254/255/257) up by 4, which
# lines up exactly with the target rows 250/251/253.
$ws.Rows("249:252").Delete() | Out-Null

# Shrink the print area to match the shorter sheet.
$ws.PageSetup.PrintArea = '$A$1:$M$248'

# Reselect / rescroll to where the deleted rows used to be, matching the
# author's saved cursor position.
$ws.Range("A249:XFD252").Select() | Out-Null
